# Update "想去人数" (F column) counts on both the "展览" and "全部类型"
# sheets, which hold duplicate data for the convention listing.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    2  = 266
    4  = 158
    7  = 102
    11 = 4620
    12 = 6884
    17 = 56
    19 = 662
    31 = 1632
    34 = 296
    40 = 156
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
